# Updates the cryptos list (Price / Volume(1h) columns) with refreshed
# values from the latest scrape, matching the GitHub Actions commit.
# A few coins (Stellar/HuobiToken, TheSandbox/Hedera) also swapped rank
# positions, so their Coin/Link columns are updated too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.598.41"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.753.82"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.45"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4498"
$ws.Range("E7").Value = "  +6.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3592"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07495"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.01"
$ws.Range("E10").Value = "  -5.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.095"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.105"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.753.27"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.37"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.826"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.649.07"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.098"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.63"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.954.28"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.075"
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.69"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("E31").Value = "  -7.79%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09094"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.666"
$ws.Range("E33").Value = "  +4.76%  "
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.00"
$ws.Range("E35").Value = "  -5.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02292"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2104"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6386"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06027"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.202"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.778"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.29"
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5921"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.953"
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.145"
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  -2.29%  "
